# Generate Report for Archive
#
# The localization status report is regenerated:
#  - Every cell whose status is "Ready for handoff" is updated to "In Translation"
#    (Overview sheet's per-locale status columns, and the "Status" column on the
#    zh-cn and de-de detail sheets).
#  - The now-narrower "Status" columns are resized to match the shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E & F) for both data rows ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn detail sheet: "Status" column (C) for both data rows ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de detail sheet: "Status" column (C) for both data rows ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Resize the status columns now that the text is shorter ---
# Target stored column width is ~13.41 chars; ColumnWidth accounts for the
# standard 5/6-character cell padding, so subtract that before assigning.
$newColumnWidth = 13.4101845877511 - (5 / 6)

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
